# Update the "Enterprises density (per 1000 people)" row for SMEs (C11)
# and MSMEs (D11) from 1.8 -> 1.78 and 11.7 -> 11.68.
#
# These cells hold their values as text (shared strings), not numbers,
# in the source workbook. Force a Text number format before assigning
# the new value so the cell keeps its original "text" type instead of
# Excel auto-converting the numeric-looking string into a real number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "1.78"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "11.68"
